# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve-profit tracking sheets
# (columns H..N: currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 88485.35000000001
$ws.Cells.Item(9, 9).Value = 107380.57
$ws.Cells.Item(9, 10).Value = 307.66666
$ws.Cells.Item(9, 11).Value = 107380.57
$ws.Cells.Item(9, 12).Value = 307.66666
$ws.Cells.Item(9, 13).Value = -107211.57
$ws.Cells.Item(9, 14).Value = -645.66666
$ws.Cells.Item(19, 8).Value = 24075652
$ws.Cells.Item(19, 9).Value = 41667856
$ws.Cells.Item(19, 10).Value = 10001889
$ws.Cells.Item(19, 11).Value = 41667856
$ws.Cells.Item(19, 12).Value = 10001889
$ws.Cells.Item(19, 13).Value = -41667681
$ws.Cells.Item(19, 14).Value = -10002239
$ws.Cells.Item(41, 8).Value = 1971.3334
$ws.Cells.Item(41, 10).Value = 930.3333
$ws.Cells.Item(41, 12).Value = 930.3333
$ws.Cells.Item(41, 14).Value = -1810.3333
$ws.Cells.Item(51, 8).Value = 7847.619
$ws.Cells.Item(51, 9).Value = 50300
$ws.Cells.Item(51, 10).Value = 5725
$ws.Cells.Item(51, 11).Value = 50300
$ws.Cells.Item(51, 12).Value = 5725
$ws.Cells.Item(51, 13).Value = -49816
$ws.Cells.Item(51, 14).Value = -6693
$ws.Cells.Item(76, 8).Value = 10899.5
$ws.Cells.Item(76, 9).Value = 18332.334
$ws.Cells.Item(76, 10).Value = 3466.6667
$ws.Cells.Item(76, 11).Value = 18332.334
$ws.Cells.Item(76, 12).Value = 3466.6667
$ws.Cells.Item(76, 13).Value = -18017.334
$ws.Cells.Item(76, 14).Value = -4096.6667
$ws.Cells.Item(79, 8).Value = 10899.5
$ws.Cells.Item(79, 9).Value = 18332.334
$ws.Cells.Item(79, 10).Value = 3466.6667
$ws.Cells.Item(79, 11).Value = 18332.334
$ws.Cells.Item(79, 12).Value = 3466.6667
$ws.Cells.Item(79, 13).Value = -17240.334
$ws.Cells.Item(79, 14).Value = -5650.6667
$ws.Cells.Item(96, 8).Value = 691889.1
$ws.Cells.Item(96, 9).Value = 1855.6364
$ws.Cells.Item(96, 10).Value = 1450925.9
$ws.Cells.Item(96, 11).Value = 5566.9092
$ws.Cells.Item(96, 12).Value = 4352777.699999999
$ws.Cells.Item(96, 13).Value = -4193.9092
$ws.Cells.Item(96, 14).Value = -4355523.699999999
$ws.Cells.Item(130, 8).Value = 82499.5
$ws.Cells.Item(130, 10).Value = 82499.5
$ws.Cells.Item(130, 12).Value = 82499.5
$ws.Cells.Item(130, 14).Value = -92539.5
$ws.Cells.Item(132, 8).Value = 5729.9062
$ws.Cells.Item(132, 9).Value = 2617.3333
$ws.Cells.Item(132, 11).Value = 7851.999899999999
$ws.Cells.Item(132, 13).Value = -5321.999899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5339
$ws.Cells.Item(2, 9).Value = 7003
$ws.Cells.Item(2, 11).Value = 7003
$ws.Cells.Item(2, 13).Value = -6890
$ws.Cells.Item(63, 8).Value = 2794.4
$ws.Cells.Item(63, 9).Value = 2824
$ws.Cells.Item(63, 10).Value = 2750
$ws.Cells.Item(63, 11).Value = 2824
$ws.Cells.Item(63, 12).Value = 2750
$ws.Cells.Item(63, 13).Value = -2138
$ws.Cells.Item(63, 14).Value = -4122
$ws.Cells.Item(66, 8).Value = 2794.4
$ws.Cells.Item(66, 9).Value = 2824
$ws.Cells.Item(66, 10).Value = 2750
$ws.Cells.Item(66, 11).Value = 14120
$ws.Cells.Item(66, 12).Value = 13750
$ws.Cells.Item(66, 13).Value = -10688
$ws.Cells.Item(66, 14).Value = -20614
$ws.Cells.Item(116, 8).Value = 5339
$ws.Cells.Item(116, 9).Value = 7003
$ws.Cells.Item(116, 11).Value = 7003
$ws.Cells.Item(116, 13).Value = -4709
$ws.Cells.Item(132, 8).Value = 1697758.9
$ws.Cells.Item(132, 9).Value = 2688.5818
$ws.Cells.Item(132, 11).Value = 8065.7454
$ws.Cells.Item(132, 13).Value = -5535.7454

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5339
$ws.Cells.Item(3, 9).Value = 7003
$ws.Cells.Item(3, 11).Value = 7003
$ws.Cells.Item(3, 13).Value = -6889
$ws.Cells.Item(99, 8).Value = 2976.6155
$ws.Cells.Item(99, 9).Value = 2539.8
$ws.Cells.Item(99, 10).Value = 4432.6665
$ws.Cells.Item(99, 11).Value = 2539.8
$ws.Cells.Item(99, 12).Value = 4432.6665
$ws.Cells.Item(99, 13).Value = -1041.8
$ws.Cells.Item(99, 14).Value = -7428.6665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 328.7143
$ws.Cells.Item(7, 9).Value = 63.25
$ws.Cells.Item(7, 10).Value = 682.6667
$ws.Cells.Item(7, 11).Value = 63.25
$ws.Cells.Item(7, 12).Value = 682.6667
$ws.Cells.Item(7, 13).Value = 49.75
$ws.Cells.Item(7, 14).Value = -908.6667
$ws.Cells.Item(31, 8).Value = 43481104
$ws.Cells.Item(31, 9).Value = 55558308
$ws.Cells.Item(31, 10).Value = 3158.2
$ws.Cells.Item(31, 11).Value = 55558308
$ws.Cells.Item(31, 12).Value = 3158.2
$ws.Cells.Item(31, 13).Value = -55558013
$ws.Cells.Item(31, 14).Value = -3748.2
$ws.Cells.Item(34, 8).Value = 43481104
$ws.Cells.Item(34, 9).Value = 55558308
$ws.Cells.Item(34, 10).Value = 3158.2
$ws.Cells.Item(34, 11).Value = 55558308
$ws.Cells.Item(34, 12).Value = 3158.2
$ws.Cells.Item(34, 13).Value = -55558106
$ws.Cells.Item(34, 14).Value = -3562.2
$ws.Cells.Item(99, 8).Value = 41465.11
$ws.Cells.Item(99, 9).Value = 11883.714
$ws.Cells.Item(99, 10).Value = 145000
$ws.Cells.Item(99, 11).Value = 11883.714
$ws.Cells.Item(99, 12).Value = 145000
$ws.Cells.Item(99, 13).Value = -10385.714
$ws.Cells.Item(99, 14).Value = -147996
$ws.Cells.Item(107, 8).Value = 5057.25
$ws.Cells.Item(107, 9).Value = 4100
$ws.Cells.Item(107, 11).Value = 4100
$ws.Cells.Item(107, 13).Value = -2180
$ws.Cells.Item(122, 8).Value = 3762.7144
$ws.Cells.Item(122, 9).Value = 3578.2222
$ws.Cells.Item(122, 10).Value = 4094.8
$ws.Cells.Item(122, 11).Value = 10734.6666
$ws.Cells.Item(122, 12).Value = 12284.4
$ws.Cells.Item(122, 13).Value = -8284.6666
$ws.Cells.Item(122, 14).Value = -17184.4
$ws.Cells.Item(126, 8).Value = 41465.11
$ws.Cells.Item(126, 9).Value = 11883.714
$ws.Cells.Item(126, 10).Value = 145000
$ws.Cells.Item(126, 11).Value = 35651.142
$ws.Cells.Item(126, 12).Value = 435000
$ws.Cells.Item(126, 13).Value = -33181.142
$ws.Cells.Item(126, 14).Value = -439940
$ws.Cells.Item(134, 8).Value = 2204.2727
$ws.Cells.Item(134, 9).Value = 1999.7
$ws.Cells.Item(134, 11).Value = 5999.1
$ws.Cells.Item(134, 13).Value = -3464.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 76923240
$ws.Cells.Item(12, 9).Value = 333333470
$ws.Cells.Item(12, 10).Value = 167.2
$ws.Cells.Item(12, 11).Value = 1000000410
$ws.Cells.Item(12, 12).Value = 501.6
$ws.Cells.Item(12, 13).Value = -1000000237
$ws.Cells.Item(12, 14).Value = -847.5999999999999
$ws.Cells.Item(136, 8).Value = 6038.4165
$ws.Cells.Item(136, 9).Value = 1239.2222
$ws.Cells.Item(136, 11).Value = 3717.6666
$ws.Cells.Item(136, 13).Value = 1382.3334
$ws.Cells.Item(138, 8).Value = 10581.857
$ws.Cells.Item(138, 9).Value = 9457.5
$ws.Cells.Item(138, 10).Value = 17328
$ws.Cells.Item(138, 11).Value = 28372.5
$ws.Cells.Item(138, 12).Value = 51984
$ws.Cells.Item(138, 13).Value = -23232.5
$ws.Cells.Item(138, 14).Value = -62264

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3004.75
$ws.Cells.Item(102, 9).Value = 3020
$ws.Cells.Item(102, 11).Value = 3020
$ws.Cells.Item(102, 13).Value = -1398
$ws.Cells.Item(113, 8).Value = 2649191.2
$ws.Cells.Item(113, 10).Value = 6175805.5
$ws.Cells.Item(113, 12).Value = 6175805.5
$ws.Cells.Item(113, 14).Value = -6180145.5
$ws.Cells.Item(122, 8).Value = 3422.9473
$ws.Cells.Item(122, 9).Value = 3296.2354
$ws.Cells.Item(122, 10).Value = 4500
$ws.Cells.Item(122, 11).Value = 9888.706200000001
$ws.Cells.Item(122, 12).Value = 13500
$ws.Cells.Item(122, 13).Value = -7438.706200000001
$ws.Cells.Item(122, 14).Value = -18400
$ws.Cells.Item(132, 8).Value = 9093446
$ws.Cells.Item(132, 9).Value = 2790.6
$ws.Cells.Item(132, 11).Value = 8371.799999999999
$ws.Cells.Item(132, 13).Value = -5841.799999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8405.272000000001
$ws.Cells.Item(7, 9).Value = 7516.2856
$ws.Cells.Item(7, 11).Value = 7516.2856
$ws.Cells.Item(7, 13).Value = -7404.2856
$ws.Cells.Item(16, 8).Value = 3531.9565
$ws.Cells.Item(16, 9).Value = 1952.25
$ws.Cells.Item(16, 11).Value = 1952.25
$ws.Cells.Item(16, 13).Value = -1782.25
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(46, 8).Value = 10000
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 10000
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 10000
$ws.Cells.Item(68, 8).Value = 5210668
$ws.Cells.Item(68, 9).Value = 6946424
$ws.Cells.Item(68, 10).Value = 3401.5
$ws.Cells.Item(68, 11).Value = 6946424
$ws.Cells.Item(68, 12).Value = 3401.5
$ws.Cells.Item(68, 13).Value = -6945675
$ws.Cells.Item(68, 14).Value = -4899.5
$ws.Cells.Item(71, 8).Value = 5210668
$ws.Cells.Item(71, 9).Value = 6946424
$ws.Cells.Item(71, 10).Value = 3401.5
$ws.Cells.Item(71, 11).Value = 34732120
$ws.Cells.Item(71, 12).Value = 17007.5
$ws.Cells.Item(71, 13).Value = -34728376
$ws.Cells.Item(71, 14).Value = -24495.5
$ws.Cells.Item(82, 8).Value = 4897.9287
$ws.Cells.Item(82, 9).Value = 2468.3
$ws.Cells.Item(82, 10).Value = 10972
$ws.Cells.Item(82, 11).Value = 2468.3
$ws.Cells.Item(82, 12).Value = 10972
$ws.Cells.Item(82, 13).Value = -2107.3
$ws.Cells.Item(82, 14).Value = -11694
$ws.Cells.Item(85, 8).Value = 4897.9287
$ws.Cells.Item(85, 9).Value = 2468.3
$ws.Cells.Item(85, 10).Value = 10972
$ws.Cells.Item(85, 11).Value = 2468.3
$ws.Cells.Item(85, 12).Value = 10972
$ws.Cells.Item(85, 13).Value = -1220.3
$ws.Cells.Item(85, 14).Value = -13468
$ws.Cells.Item(93, 8).Value = 1545549.2
$ws.Cells.Item(93, 9).Value = 774
$ws.Cells.Item(93, 10).Value = 9269426
$ws.Cells.Item(93, 11).Value = 774
$ws.Cells.Item(93, 12).Value = 9269426
$ws.Cells.Item(93, 13).Value = 474
$ws.Cells.Item(93, 14).Value = -9271922
$ws.Cells.Item(126, 8).Value = 8405.272000000001
$ws.Cells.Item(126, 9).Value = 7516.2856
$ws.Cells.Item(126, 11).Value = 22548.8568
$ws.Cells.Item(126, 13).Value = -20078.8568
$ws.Cells.Item(132, 8).Value = 3037.6206
$ws.Cells.Item(132, 9).Value = 2003.9166
$ws.Cells.Item(132, 11).Value = 6011.7498
$ws.Cells.Item(132, 13).Value = -3481.7498
$ws.Cells.Item(40, 13).ClearContents()  # remove M40 (was -1331.3334)
$ws.Cells.Item(46, 13).ClearContents()  # remove M46 (was -1799.5)
$ws.Cells.Item(46, 14).Value = -10376  # add N46

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3137.6667
$ws.Cells.Item(122, 9).Value = 2294.8462
$ws.Cells.Item(122, 10).Value = 4133.727
$ws.Cells.Item(122, 11).Value = 6884.5386
$ws.Cells.Item(122, 12).Value = 12401.181
$ws.Cells.Item(122, 13).Value = -4434.5386
$ws.Cells.Item(122, 14).Value = -17301.181
$ws.Cells.Item(136, 8).Value = 288328.66
$ws.Cells.Item(136, 9).Value = 2390.2903
$ws.Cells.Item(136, 10).Value = 2504351
$ws.Cells.Item(136, 11).Value = 7170.8709
$ws.Cells.Item(136, 12).Value = 7513053
$ws.Cells.Item(136, 13).Value = -4620.8709
$ws.Cells.Item(136, 14).Value = -7518153
